$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173, shifting existing rows down
$ws.Rows.Item(173).Insert()

# Populate the new row with the "finishbet" entry
$ws.Range("A173").Value2 = "finishbet"
$ws.Range("B173").Value2 = "Maks 100k çekim"
$ws.Range("C173").Value2 = "yatırımsız"
$ws.Range("D173").Value2 = "Evet"

# Restore the selection to B13 (matches the saved view state in the target file)
$ws.Range("B13").Select()
